$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C ("Förändrad") from 45233 to 45243 for rows 2 through 18
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2() -eq 45233) {
        $cell.Value = 45243
    }
}
